$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.480879068374634
$ws.Range("B1").Value = 1.999357223510742
$ws.Range("C1").Value = 2.301314830780029
$ws.Range("D1").Value = 2.761540412902832
$ws.Range("E1").Value = 2.750622510910034
